$d = $word.ActiveDocument

# 1. Remove the leftover "_GoBack" bookmark that Word drops after the
#    cursor's last edit location (after "...scores.").
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    if ($goBack -ne $null) {
        $goBack.Delete()
    }
} catch {
}

# 2. "Please use the contact page in case of questions or issues." ->
#    "Please use the contact page for questions or issues."
#    (both occurrences in the document)
$rng = $d.Content
[void]$rng.Find.Execute("in case of", $true, $false, $false, $false, $false, $true, 1, $false, "for", 2)

# 3. Nudge the second "click" screenshot image slightly (anchor offset
#    changed from (830925, 17145) EMU to (830580, 21156) EMU).
$shape = $d.Shapes.Item(2)
$shape.Left = 65.4
$shape.Top = 1.6658267716535433
